$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.032.98"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -3.17%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.714.61"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -3.00%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.05%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'310.02"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -5.65%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.08%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4753"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +4.47%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.3448"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -2.24%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'42.21"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.49%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.07242"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -1.90%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  -5.06%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'1.002"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.13%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'19.75"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -4.65%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'5.820"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -3.16%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'1.717.13"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -2.94%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'6.808"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -5.23%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'87.12"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -5.85%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  -2.41%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'0.06371"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -1.08%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.08%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'16.42"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -3.15%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'5.606"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -2.62%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'27.094.80"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -3.04%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'10.71"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -4.33%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.096"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.18%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'19.96"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.98%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'150.90"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -4.78%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'1.916.66"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -2.88%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  -3.94%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'120.74"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -2.52%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'1.022"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -5.45%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.09212"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.32%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'3.602"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -1.71%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'5.311"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -5.31%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'1.472"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +6.37%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.02180"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -4.29%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.05832"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -4.59%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'10.92"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -7.72%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.1989"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -4.88%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'1.001"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.11%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'4.705"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -4.90%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.5952"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -4.83%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'1.084"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -8.13%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'7.490"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -3.90%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'12.79"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -3.51%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'3.581"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -4.35%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.5563"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -4.92%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'118.82"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.99%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'1.824"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -5.39%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'1.106"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -2.20%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.06635"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -2.85%  "
$ws.Range('E51').Style = 'Normal'
